$wb = $excel.ActiveWorkbook

# --- Update the product name text (a hyphen was inserted after "197") ---
# This shared string is used by both sheets' B1 cell.
$oldName = "197MS-EI-DB-DL-REC-NON-RNI-CTPD-SAR-MD-TR-1-ONTIME"
$newName = "197-MS-EI-DB-DL-REC-NON-RNI-CTPD-SAR-MD-TR-1-ONTIME"

$wsInput = $wb.Worksheets.Item("ProductLoanInput")
$wsOutput = $wb.Worksheets.Item("ProductLoanOutput")

$wsInput.Range("B1").Value = $newName
$wsOutput.Range("B1").Value = $newName

# --- Update view/selection state ---
# Select B1 on the input sheet (it is being left as the non-active tab).
[void]$wsInput.Activate()
[void]$wsInput.Range("B1").Select()

# Make the output sheet the active tab, with B1 selected there.
[void]$wsOutput.Activate()
[void]$wsOutput.Range("B1").Select()
